# Travis County 2020 bg SVI workbook update
# Re-run of the factor-analysis pipeline: factor member ordering changed
# and loading/variance values refreshed with newly computed (numerically
# close) results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Significant Components": reorder the members listed per factor
# ---------------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significant Components")
$wsSig.Range("C2").Value = "['QEXTRCT' 'QMOHO' 'QESL' 'QHISPC' 'QEDLESHI' 'PPUNIT' 'QNOHLTH']"
$wsSig.Range("C3").Value = "['QRICH' 'PERCAP' 'MDHSEVAL']"
$wsSig.Range("C4").Value = "['QAGEDEP' 'MEDAGE' 'QSSBEN']"
$wsSig.Range("C6").Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"

# ---------------------------------------------------------------------
# Sheet "Loading Factors": rows reordered + refreshed values
# ---------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Loading Factors")

$loadData = @(
    @("QEXTRCT",  0.7662017872076399,  0.1246192454280793,  0.05703177946023106,  0.1061640317557724,  -0.1262460487065369,  0.004275829328621803),
    @("QMOHO",    0.4515387057450205,  0.1342118383062358,  0.03533061706842234,  -0.1630772611730658, -0.01675659515438632,  0.001862667846193435),
    @("QESL",     0.7463154807774712,  0.1383746659566736,  -0.06828868611331276, 0.2025402326886757,  -0.0440500951558544,   0.09086865729578704),
    @("QHISPC",   0.7423692411176009,  0.3098526185856265,  -0.14869528158947,    0.09400583920772702, -0.04885413337778412,  0.2354826915134233),
    @("QEDLESHI", 0.835277940995633,   0.1442623932411957,  0.002317829985271983, 0.1530479583262567,  -0.02604026930732095,  0.1765466364949183),
    @("PPUNIT",   0.5835172974385709,  -0.03689746289890063,-0.1272499891103111,  -0.4204908261492347, 0.1006093493267103,    0.1309775141776323),
    @("QNOHLTH",  0.6038940737830031,  0.3577588836664262,  -0.07977624915877934, 0.2571839839282107,  -0.09534314929099239,  0.1820152844613006),
    @("QRICH",    0.1786760057396259,  0.8099499507036548,  -0.11711372325679,    0.368662757418611,   -0.02002121838050307,  0.1095116346276093),
    @("PERCAP",   0.4141283906093469,  0.6768765560650708,  -0.2089607222700036,  0.1891635479774894,  0.09876520938485223,   0.1797756318032703),
    @("MDHSEVAL", 0.3537576154528908,  0.7287084011579007,  -0.08538345335431308, -0.04214449815258696,0.04384743061208778,   0.1201375913882453),
    @("QAGEDEP",  0.03196500340378081, -0.06645709741022643,0.6956263865925336,   -0.04457261749977444,0.6123584626070012,    -0.04113476760492205),
    @("QRENTER",  -0.06964049953173979,0.2243052263311398,  -0.4250942600076245,  0.7486728956473566,  -0.08020597001281483,  0.06481864854065772),
    @("MEDAGE",   -0.1933922649374947, -0.255600559253691,  0.773581617477663,    -0.2675924857923357, -0.05058706149342817,  -0.1226097599073551),
    @("QSSBEN",   0.02114278564884536, -0.03947992304733914,0.7957449321843202,   -0.08782706832591999,0.06483514891930618,   0.04245536416394331),
    @("QNOAUTO",  0.05709370143238238, 0.04107575801101934, 0.01144308936803343,  0.6493597835740099,  0.01494328553999053,   0.05784573716820147),
    @("QPOVTY",   0.1957216577902933,  0.1330977694032092,  -0.2320824168709502,  0.5572080780311643,  0.08212745642362895,   0.1166838083140116),
    @("QFEMLBR",  -0.2164852247003474, 0.07008131208836509, -0.05054873831355557, -0.01544327866450019,0.6764703221063801,    0.13279831781132),
    @("QFEMALE",  0.02896968261562555, 0.001264419460793409,0.1484947708258256,   0.04764862546510356, 0.9303837721033712,    -0.0001998282843647516),
    @("QFHH",     0.3059368149044416,  0.1540730305838567,  0.02452183686754002,  0.02957443198154009, 0.1777865093132783,    0.9169224231249294),
    @("QFAM",     0.1389414206732875,  0.2227513374345099,  -0.1031614094357772,  0.3534990165047647,  -0.02460918477747409,  0.4969878297920926)
)

for ($i = 0; $i -lt $loadData.Length; $i++) {
    $row = $i + 2
    $rowdata = $loadData[$i]
    $wsLoad.Cells.Item($row, 1).Value = $rowdata[0]
    for ($j = 1; $j -le 6; $j++) {
        $wsLoad.Cells.Item($row, $j + 1).Value = $rowdata[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet "All Refactor Variances": refresh the F1 columns (I:N), rows 2-5
# ---------------------------------------------------------------------
$wsAllVar = $wb.Worksheets.Item("All Refactor Variances")

$allVarData = @(
    @(3.876385837741227, 2.163351746552688, 2.104589608430914, 2.029886901840213, 1.803101118894937, 1.348998911800272),
    @(0.1938192918870613, 0.1081675873276344, 0.1052294804215457, 0.1014943450920106, 0.09015505594474683, 0.06744994559001359),
    @(0.1938192918870613, 0.3019868792146957, 0.4072163596362415, 0.5087107047282521, 0.5988657606729989, 0.6663157062630125),
    @(0.2908820699636273, 0.1623368416966863, 0.1579273600073429, 0.1523217059691343, 0.1353038133385381, 0.1012282090246711)
)

for ($i = 0; $i -lt $allVarData.Length; $i++) {
    $row = $i + 2
    $rowdata = $allVarData[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $wsAllVar.Cells.Item($row, $j + 9).Value = $rowdata[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet "Final Variances": refresh B2:G5 (mirrors the F1 columns above)
# ---------------------------------------------------------------------
$wsFinalVar = $wb.Worksheets.Item("Final Variances")

$finalVarData = @(
    @(3.876385837741227, 2.163351746552688, 2.104589608430914, 2.029886901840213, 1.803101118894937, 1.348998911800272),
    @(0.1938192918870613, 0.1081675873276344, 0.1052294804215457, 0.1014943450920106, 0.09015505594474683, 0.06744994559001359),
    @(0.1938192918870613, 0.3019868792146957, 0.4072163596362415, 0.5087107047282521, 0.5988657606729989, 0.6663157062630125),
    @(0.2908820699636273, 0.1623368416966863, 0.1579273600073429, 0.1523217059691343, 0.1353038133385381, 0.1012282090246711)
)

for ($i = 0; $i -lt $finalVarData.Length; $i++) {
    $row = $i + 2
    $rowdata = $finalVarData[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $wsFinalVar.Cells.Item($row, $j + 2).Value = $rowdata[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet "Included and Excluded": reorder the included-variable list
# ---------------------------------------------------------------------
$wsIncExc = $wb.Worksheets.Item("Included and Excluded")
$wsIncExc.Range("B2").Value = "[['QEXTRCT', 'QMOHO', 'QESL', 'QHISPC', 'QEDLESHI', 'PPUNIT', 'QNOHLTH', 'QRICH', 'PERCAP', 'MDHSEVAL', 'QAGEDEP', 'MEDAGE', 'QSSBEN', 'QRENTER', 'QNOAUTO', 'QPOVTY', 'QFEMLBR', 'QFEMALE', 'QFHH', 'QFAM']]"
